# Updated cryptos list on Sun Apr 30 18:31:54 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" '29.847.39'
Set-TextValue "D3" '1.936.05'
Set-TextValue "E3" '  +1.21%  '
Set-TextValue "D4" '1.008'
Set-TextValue "E4" '  -0.20%  '
Set-TextValue "D5" '336.26'
Set-TextValue "E5" '  +3.31%  '
Set-TextValue "D6" '1.007'
Set-TextValue "E6" '  -0.16%  '
Set-TextValue "D7" '0.4836'
Set-TextValue "E7" '  +0.47%  '
Set-TextValue "D8" '0.4107'
Set-TextValue "E8" '  +1.04%  '
Set-TextValue "E9" '  -0.86%  '
Set-TextValue "D10" '1.014'
Set-TextValue "E10" '  -0.66%  '
Set-TextValue "D11" '23.70'
Set-TextValue "E11" '  +0.99%  '
Set-TextValue "D12" '1.968.50'
Set-TextValue "E12" '  +3.62%  '
Set-TextValue "D13" '6.080'
Set-TextValue "E13" '  +0.75%  '
Set-TextValue "D14" '7.280'
Set-TextValue "E14" '  +1.02%  '
Set-TextValue "D15" '90.86'
Set-TextValue "E15" '  -0.28%  '
Set-TextValue "D16" '0.06833'
Set-TextValue "E16" '  +0.51%  '
Set-TextValue "E17" '  -0.16%  '
Set-TextValue "E18" '  -0.56%  '
Set-TextValue "D19" '17.75'
Set-TextValue "E19" '  +0.35%  '
Set-TextValue "D21" '29.828.27'
Set-TextValue "E21" '  +1.28%  '
Set-TextValue "D22" '5.623'
Set-TextValue "E22" '  +0.03%  '
Set-TextValue "E23" '  +0.30%  '
Set-TextValue "D24" '2.177'
Set-TextValue "E24" '  -0.82%  '
Set-TextValue "D25" '2.170.64'
Set-TextValue "E25" '  +3.32%  '
Set-TextValue "D26" '6.600'
Set-TextValue "E26" '  +0.21%  '
Set-TextValue "D27" '156.59'
Set-TextValue "E27" '  -0.05%  '
Set-TextValue "D28" '20.06'
Set-TextValue "E28" '  +0.28%  '
Set-TextValue "D29" '2.093'
Set-TextValue "E29" '  -0.51%  '
Set-TextValue "E30" '  +0.66%  '
Set-TextValue "D31" '1.005'
Set-TextValue "E31" '  -1.38%  '
Set-TextValue "D32" '0.09647'
Set-TextValue "E32" '  +0.93%  '
Set-TextValue "D33" '5.551'
Set-TextValue "E33" '  -0.21%  '
Set-TextValue "D34" '1.413'
Set-TextValue "E34" '  +3.48%  '
Set-TextValue "D35" '3.532'
Set-TextValue "E35" '  -0.58%  '
Set-TextValue "D36" '0.06565'
Set-TextValue "E36" '  +7.35%  '
Set-TextValue "D37" '0.02284'
Set-TextValue "E37" '  +0.02%  '
Set-TextValue "E38" '  +1.94%  '
Set-TextValue "D39" '0.5964'
Set-TextValue "E39" '  +0.04%  '
Set-TextValue "E40" '  -0.81%  '
Set-TextValue "D41" '7.945'
Set-TextValue "E41" '  -1.19%  '
Set-TextValue "D42" '0.1847'
Set-TextValue "E42" '  -0.09%  '
Set-TextValue "D43" '2.468'
Set-TextValue "E43" '  +3.12%  '
Set-TextValue "D44" '1.276'
Set-TextValue "E44" '  -0.45%  '
Set-TextValue "D45" '12.29'
Set-TextValue "E45" '  -1.34%  '
Set-TextValue "D46" '0.07476'
Set-TextValue "E46" '  -1.71%  '
Set-TextValue "D47" '0.5557'
Set-TextValue "E47" '  -0.26%  '
Set-TextValue "D48" '1.987'
Set-TextValue "E48" '  +1.75%  '
Set-TextValue "D49" '116.81'
Set-TextValue "E49" '  -0.58%  '
Set-TextValue "D50" '2.422'
Set-TextValue "E50" '  -0.30%  '
Set-TextValue "D51" '72.57'
Set-TextValue "E51" '  +0.53%  '
